$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change shared string value "E7760" -> "E7420" (column H, rows 2-32)
$ws.Range("H2:H32").Value = "E7420"

# 2. Change the selection on the sheet from E2:E32 to H2:H32
$ws.Range("H2:H32").Select()

$wb.Save()
